$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "MiddleName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "EmployeeID"

$ws.Range("D1").Select()
